# Mark the following DSA question rows as answered ("yes") in column C,
# matching the "dp and binary tree" study-progress update.
#
# Rows 156-165  : Doubly Linked List rotation/reverse questions
# Rows 177-180  : Binary Trees - traversal / height / diameter questions
# Rows 356-357,359 : DP questions
# Rows 412,416,418-419 : DP questions
# Rows 446-447  : DP questions
# Row 472       : DP questions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @(
    "C156","C157","C158","C159","C160","C161","C162","C163","C164","C165",
    "C177","C178","C179","C180",
    "C356","C357","C359",
    "C412","C416","C418","C419",
    "C446","C447",
    "C472"
)

# C358 already carries the "yes" formatting/value we want to replicate
# (green fill, style index 12) onto the target cells.
$src = $ws.Range("C358")
$src.Copy()

foreach ($addr in $cells) {
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

foreach ($addr in $cells) {
    $ws.Range($addr).Value = "yes"
}

# Reflect the editor's final cursor position/selection from the edit.
$ws.Activate()
$ws.Range("E179").Select()
